# edit.ps1
# Applies the Lgi4-Adam22 LR-pairs update: new TPM-derived values, and the
# sending-cluster block for "ECs" is dropped while the remaining two blocks
# (FAPs, MuSCs as sending clusters) are refreshed with new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 8-10 (the old "sending cluster = MuSCs" block) - data now has
# 6 data rows instead of 9, so the sheet shrinks from A1:T10 to A1:T7.
$ws.Range("A8:T10").EntireRow.Delete()

# New data for rows 2-7 (header row 1 is unchanged).
$data = @(
    @("FAPs", "Lgi4", "Adam22", "ECs", 3, 1, 10.326286, 30.978858, 0.5836270415642173, 0.5836270415642173, 3, 1, 3.618510333333333, 10.855531, 0.1815566256530994, 0.1815566256530994, 37.36577259595533, 336.2919533635979, 0.1059613563063005, 0.1059613563063005),
    @("FAPs", "Lgi4", "Adam22", "FAPs", 3, 1, 10.326286, 30.978858, 0.5836270415642173, 0.5836270415642173, 3, 1, 3.860692, 11.582076, 0.1937079481987336, 0.1937079481987336, 39.866609749912, 358.799487749208, 0.1130531967347016, 0.1130531967347016),
    @("FAPs", "Lgi4", "Adam22", "MuSCs", 3, 1, 10.326286, 30.978858, 0.5836270415642173, 0.5836270415642173, 3, 1, 12.45127566666667, 37.353827, 0.6247354261481669, 0.6247354261481669, 128.5754335988406, 1157.178902389566, 0.3646124885232152, 0.3646124885232152),
    @("MuSCs", "Lgi4", "Adam22", "ECs", 3, 1, 7.367010000000001, 22.10103, 0.4163729584357827, 0.4163729584357828, 3, 1, 3.618510333333333, 10.855531, 0.1815566256530994, 0.1815566256530994, 26.65760181077, 239.91841629693, 0.0755952693467989, 0.07559526934679893),
    @("MuSCs", "Lgi4", "Adam22", "FAPs", 3, 1, 7.367010000000001, 22.10103, 0.4163729584357827, 0.4163729584357828, 3, 1, 3.860692, 11.582076, 0.1937079481987336, 0.1937079481987336, 28.44175657092, 255.97580913828, 0.08065475146403207, 0.08065475146403207),
    @("MuSCs", "Lgi4", "Adam22", "MuSCs", 3, 1, 7.367010000000001, 22.10103, 0.4163729584357827, 0.4163729584357828, 3, 1, 12.45127566666667, 37.353827, 0.6247354261481669, 0.6247354261481669, 91.72867234909, 825.55805114181, 0.2601229376249517, 0.2601229376249518)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
